{"js": "// Replace the date line and the 25 division-problem cells with their\n// updated values, matched by exact original text (each original value\n// is unique in the document, so a scoped search-and-replace is safe).\nconst replacements = [\n  [\"2023-11-14 Tuesday\", \"2023-11-15 Wednesday\"],\n  [\"77\u00f77=\", \"65\u00f77=\"],\n  [\"10\u00f77=\", \"50\u00f77=\"],\n  [\"35\u00f76=\", \"36\u00f72=\"],\n  [\"21\u00f73=\", \"19\u00f78=\"],\n  [\"96\u00f77=\", \"10\u00f74=\"],\n  [\"86\u00f74=\", \"20\u00f72=\"],\n  [\"80\u00f74=\", \"35\u00f74=\"],\n  [\"93\u00f77=\", \"41\u00f78=\"],\n  [\"56\u00f73=\", \"90\u00f72=\"],\n  [\"21\u00f75=\", \"32\u00f77=\"],\n  [\"28\u00f78=\", \"32\u00f74=\"],\n  [\"98\u00f79=\", \"40\u00f72=\"],\n  [\"27\u00f77=\", \"14\u00f78=\"],\n  [\"20\u00f78=\", \"95\u00f76=\"],\n  [\"50\u00f76=\", \"72\u00f76=\"],\n  [\"30\u00f78=\", \"48\u00f73=\"],\n  [\"71\u00f74=\", \"29\u00f79=\"],\n  [\"73\u00f77=\", \"30\u00f74=\"],\n  [\"12\u00f75=\", \"95\u00f74=\"],\n  [\"37\u00f75=\", \"89\u00f77=\"],\n  [\"48\u00f75=\", \"72\u00f79=\"],\n  [\"49\u00f74=\", \"95\u00f74=\"],\n  [\"50\u00f74=\", \"48\u00f77=\"],\n  [\"44\u00f76=\", \"43\u00f79=\"],\n  [\"87\u00f78=\", \"59\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 division-problem cells with their\n# updated values, matched by exact original text (each original value\n# is unique in the document, so Find/Replace is unambiguous).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2023-11-14 Tuesday\", \"2023-11-15 Wednesday\"),\n  @(\"77\u00f77=\", \"65\u00f77=\"),\n  @(\"10\u00f77=\", \"50\u00f77=\"),\n  @(\"35\u00f76=\", \"36\u00f72=\"),\n  @(\"21\u00f73=\", \"19\u00f78=\"),\n  @(\"96\u00f77=\", \"10\u00f74=\"),\n  @(\"86\u00f74=\", \"20\u00f72=\"),\n  @(\"80\u00f74=\", \"35\u00f74=\"),\n  @(\"93\u00f77=\", \"41\u00f78=\"),\n  @(\"56\u00f73=\", \"90\u00f72=\"),\n  @(\"21\u00f75=\", \"32\u00f77=\"),\n  @(\"28\u00f78=\", \"32\u00f74=\"),\n  @(\"98\u00f79=\", \"40\u00f72=\"),\n  @(\"27\u00f77=\", \"14\u00f78=\"),\n  @(\"20\u00f78=\", \"95\u00f76=\"),\n  @(\"50\u00f76=\", \"72\u00f76=\"),\n  @(\"30\u00f78=\", \"48\u00f73=\"),\n  @(\"71\u00f74=\", \"29\u00f79=\"),\n  @(\"73\u00f77=\", \"30\u00f74=\"),\n  @(\"12\u00f75=\", \"95\u00f74=\"),\n  @(\"37\u00f75=\", \"89\u00f77=\"),\n  @(\"48\u00f75=\", \"72\u00f79=\"),\n  @(\"49\u00f74=\", \"95\u00f74=\"),\n  @(\"50\u00f74=\", \"48\u00f77=\"),\n  @(\"44\u00f76=\", \"43\u00f79=\"),\n  @(\"87\u00f78=\", \"59\u00f78=\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
